# "Agregadas ventanas de carrito y retroceso/vaciado"
# Pre-allocate 15 extra blank rows (11-25) below the existing product list so
# the cart / undo-clear feature always has room to grow into. Column A keeps
# the running Id counter; columns B:E are emptied out (ready to be filled in
# later by the app) but still present as real (text/empty) cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 11; $i -le 25; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1

    # A bare "" assignment is indistinguishable from "no cell" to this
    # engine, so force a text cell via the classic leading-apostrophe
    # "treat as text" prefix, then strip the quote-prefix styling it adds
    # so the cell ends up a plain, unstyled empty text cell.
    $ws.Cells.Item($i, 2).Value = "'"
    $ws.Cells.Item($i, 2).Style = "Normal"
    $ws.Cells.Item($i, 3).Value = "'"
    $ws.Cells.Item($i, 3).Style = "Normal"
    $ws.Cells.Item($i, 4).Value = "'"
    $ws.Cells.Item($i, 4).Style = "Normal"
    $ws.Cells.Item($i, 5).Value = "'"
    $ws.Cells.Item($i, 5).Style = "Normal"
}
